# Add the "2022-Q3" quarterly fund-holdings sheet and update the "总计"
# (totals) summary sheet to include it, matching the commit's target state:
#   总计, 2022-Q3, 2022-Q2, 2022-Q1, 2021-Q4

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Create the new "2022-Q3" sheet positioned right after "总计"
#    (i.e. right before the existing "2022-Q2" sheet).
# ---------------------------------------------------------------------------
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$q3Sheet = $wb.Worksheets.Add($q2Sheet)
$q3Sheet.Name = "2022-Q3"

# Header row (row 1) - same headers/style as the other quarterly sheets.
$q3Sheet.Range("B1").Value = "基金代码"
$q3Sheet.Range("C1").Value = "基金名称"
$q3Sheet.Range("D1").Value = "基金规模"
$q3Sheet.Range("E1").Value = "股票总仓位"
$q3Sheet.Range("F1").Value = "仓位占比"
$q3Sheet.Range("G1").Value = "持有市值(亿元)"
$q3Sheet.Range("H1").Value = "仓位排名"

# Fund holdings data for 2022-Q3 (rows 2-24).
# Columns: code, name, scale, stock-position, position-ratio, market-value, rank
$q3Data = @(
    @("011479", "广发诚享混合A", "36.88", "87.12", "5.34", "1.9694", 8),
    @("011480", "广发诚享混合C", "4.20", "87.12", "5.34", "0.2243", 8),
    @("161039", "富国中证1000指数增强（LOF）A", "25.41", "84.72", "0.60", "0.1525", 8),
    @("162717", "广发再融资主题灵活配置混合A", "1.34", "90.49", "7.25", "0.0972", 6),
    @("008602", "方正富邦新兴成长混合A", "1.23", "86.03", "6.56", "0.0807", 1),
    @("013711", "广发再融资主题灵活配置混合C", "1.04", "90.49", "7.25", "0.0754", 6),
    @("013331", "富国中证1000指数增强（LOF）C", "8.53", "84.72", "0.60", "0.0512", 8),
    @("004044", "金鹰转型动力灵活配置混合", "0.60", "88.13", "4.88", "0.0293", 9),
    @("000827", "广发中证百度百发策略100指数E", "2.74", "92.46", "1.00", "0.0274", 9),
    @("006401", "先锋量化优选灵活配置混合A", "0.60", "94.55", "2.67", "0.0160", 6),
    @("000826", "广发中证百度百发策略100指数A", "0.89", "92.46", "1.00", "0.0089", 9),
    @("006402", "先锋量化优选灵活配置混合C", "0.18", "94.55", "2.67", "0.0048", 6),
    @("008603", "方正富邦新兴成长混合C", "0.03", "86.03", "6.56", "0.0020", 1),
    @("013489", "广发东财大数据精选灵活配置混合C", "0.18", "60.68", "1.03", "0.0019", 6),
    @("002802", "广发东财大数据精选灵活配置混合A", "0.15", "60.68", "1.03", "0.0015", 6),
    @("004833", "先锋聚利灵活配置混合A", "0.05", "93.58", "2.60", "0.0013", 5),
    @("004724", "先锋聚元灵活配置混合A", "0.04", "94.36", "2.50", "0.0010", 8),
    @("004725", "先锋聚元灵活配置混合C", "0.04", "94.36", "2.50", "0.0010", 8),
    @("004727", "先锋聚优灵活配置混合C", "0.03", "91.74", "2.61", "0.0008", 3),
    @("003587", "先锋精一灵活配置混合C", "0.02", "94.32", "2.59", "0.0005", 6),
    @("004834", "先锋聚利灵活配置混合C", "0.02", "93.58", "2.60", "0.0005", 5),
    @("003586", "先锋精一灵活配置混合A", "0.01", "94.32", "2.59", "0.0003", 6),
    @("004726", "先锋聚优灵活配置混合A", "0.01", "91.74", "2.61", "0.0003", 3)
)

for ($i = 0; $i -lt $q3Data.Length; $i++) {
    $r = $i + 2
    $entry = $q3Data[$i]

    $q3Sheet.Range("A" + $r).Value = $i

    # B-G are text cells in the source data (fund codes keep leading zeros,
    # numeric-looking figures are stored as text) - force text via NumberFormat
    # so Excel doesn't silently convert them to numbers.
    $q3Sheet.Range("B" + $r + ":G" + $r).NumberFormat = "@"
    $q3Sheet.Range("B" + $r).Value = $entry[0]
    $q3Sheet.Range("C" + $r).Value = $entry[1]
    $q3Sheet.Range("D" + $r).Value = $entry[2]
    $q3Sheet.Range("E" + $r).Value = $entry[3]
    $q3Sheet.Range("F" + $r).Value = $entry[4]
    $q3Sheet.Range("G" + $r).Value = $entry[5]

    $q3Sheet.Range("H" + $r).Value = $entry[6]
}

# Drop the helper "@" number formats we used to preserve text (leaves cells
# with no explicit style, matching the plain inlineStr cells used elsewhere).
$q3Sheet.Range("B2:G24").ClearFormats()

# Copy over the header-row and index-column styling from the existing
# "2022-Q2" sheet so the new sheet matches the established look.
$q2Sheet.Range("B1:H1").Copy()
$q3Sheet.Range("B1:H1").PasteSpecial(-4122)
$q2Sheet.Range("A2").Copy()
$q3Sheet.Range("A2:A24").PasteSpecial(-4122)

$q3Sheet.Range("A1").Select()

# ---------------------------------------------------------------------------
# 2) Update the "总计" (totals) sheet: insert a new row for 2022-Q3 at the
#    top of the data (row 2), pushing the other quarters down.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 23
$totalSheet.Range("D2").Value = 2.75

# The pre-existing quarters (now rows 3-5) keep their date/count/value but
# get re-numbered in column A - index 0-based from the new top row.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3

# The inserted row inherited stray formatting - clear it and restore the
# plain look used by the other data rows (bold/centered index in col A,
# no special formatting for B-D).
$totalSheet.Range("A2").ClearFormats()
$totalSheet.Range("B2:D2").ClearFormats()
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A1").Select()
